$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 10-13 (old ECs-target rows for MuSCs/Resolving-Mac senders no longer needed,
# and the whole "ECs" target-cluster column has been dropped from the output)
$ws.Range("A10:T13").Delete() | Out-Null

$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 2.18083
$ws.Range("H2").Value = 6.54249
$ws.Range("I2").Value = 0.2029465558466658
$ws.Range("J2").Value = 0.2029465558466658
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.079670666666667
$ws.Range("N2").Value = 3.239012
$ws.Range("O2").Value = 0.5935332782376214
$ws.Range("P2").Value = 0.5935332782376214
$ws.Range("Q2").Value = 2.354578179986667
$ws.Range("R2").Value = 21.19120361988
$ws.Range("S2").Value = 0.120455534598706
$ws.Range("T2").Value = 0.120455534598706

$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 2.18083
$ws.Range("H3").Value = 6.54249
$ws.Range("I3").Value = 0.2029465558466658
$ws.Range("J3").Value = 0.2029465558466658
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.739386
$ws.Range("N3").Value = 2.218158
$ws.Range("O3").Value = 0.4064667217623787
$ws.Range("P3").Value = 0.4064667217623787
$ws.Range("Q3").Value = 1.61247517038
$ws.Range("R3").Value = 14.51227653342
$ws.Range("S3").Value = 0.08249102124795975
$ws.Range("T3").Value = 0.08249102124795975

$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 3.003971
$ws.Range("H4").Value = 9.011913
$ws.Range("I4").Value = 0.2795474971975186
$ws.Range("J4").Value = 0.2795474971975186
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.079670666666667
$ws.Range("N4").Value = 3.239012
$ws.Range("O4").Value = 0.5935332782376214
$ws.Range("P4").Value = 0.5935332782376214
$ws.Range("Q4").Value = 3.243299372217333
$ws.Range("R4").Value = 29.189694349956
$ws.Range("S4").Value = 0.1659207424347655
$ws.Range("T4").Value = 0.1659207424347655

$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 3.003971
$ws.Range("H5").Value = 9.011913
$ws.Range("I5").Value = 0.2795474971975186
$ws.Range("J5").Value = 0.2795474971975186
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.739386
$ws.Range("N5").Value = 2.218158
$ws.Range("O5").Value = 0.4064667217623787
$ws.Range("P5").Value = 0.4064667217623787
$ws.Range("Q5").Value = 2.221094101806
$ws.Range("R5").Value = 19.989846916254
$ws.Range("S5").Value = 0.1136267547627531
$ws.Range("T5").Value = 0.1136267547627531

$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("G6").Value = 3.585492333333333
$ws.Range("H6").Value = 10.756477
$ws.Range("I6").Value = 0.3336634767793112
$ws.Range("J6").Value = 0.3336634767793112
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.079670666666667
$ws.Range("N6").Value = 3.239012
$ws.Range("O6").Value = 0.5935332782376214
$ws.Range("P6").Value = 0.5935332782376214
$ws.Range("Q6").Value = 3.871150897858222
$ws.Range("R6").Value = 34.840358080724
$ws.Range("S6").Value = 0.198040377200987
$ws.Range("T6").Value = 0.198040377200987

$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("G7").Value = 3.585492333333333
$ws.Range("H7").Value = 10.756477
$ws.Range("I7").Value = 0.3336634767793112
$ws.Range("J7").Value = 0.3336634767793112
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.739386
$ws.Range("N7").Value = 2.218158
$ws.Range("O7").Value = 0.4064667217623787
$ws.Range("P7").Value = 0.4064667217623787
$ws.Range("Q7").Value = 2.651062834374
$ws.Range("R7").Value = 23.859565509366
$ws.Range("S7").Value = 0.1356230995783242
$ws.Range("T7").Value = 0.1356230995783242

$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("D8").Value = "FAPs"
$ws.Range("G8").Value = 1.975540666666667
$ws.Range("H8").Value = 5.926622
$ws.Range("I8").Value = 0.1838424701765043
$ws.Range("J8").Value = 0.1838424701765044
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.079670666666667
$ws.Range("N8").Value = 3.239012
$ws.Range("O8").Value = 0.5935332782376214
$ws.Range("P8").Value = 0.5935332782376214
$ws.Range("Q8").Value = 2.132933308607111
$ws.Range("R8").Value = 19.196399777464
$ws.Range("S8").Value = 0.1091166240031628
$ws.Range("T8").Value = 0.1091166240031628

$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = 1.975540666666667
$ws.Range("H9").Value = 5.926622
$ws.Range("I9").Value = 0.1838424701765043
$ws.Range("J9").Value = 0.1838424701765044
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.739386
$ws.Range("N9").Value = 2.218158
$ws.Range("O9").Value = 0.4064667217623787
$ws.Range("P9").Value = 0.4064667217623787
$ws.Range("Q9").Value = 1.460687111364
$ws.Range("R9").Value = 13.146184002276
$ws.Range("S9").Value = 0.07472584617334159
$ws.Range("T9").Value = 0.0747258461733416
